# Update the "Periodo Mora" period-code column (E16:E55) so that the
# periods are listed in ascending order (1612 .. 2003) instead of the
# previous descending order (2003 .. 1612). This mirrors the reordering
# of the corresponding shared-string table entries in the committed
# workbook, while keeping every other cell (value, style, border, etc.)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "1612","1701","1702","1703","1704","1705","1706","1707","1708","1709",
    "1710","1711","1712","1801","1802","1803","1804","1805","1806","1807",
    "1808","1809","1810","1811","1812","1901","1902","1903","1904","1905",
    "1906","1907","1908","1909","1910","1911","1912","2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $ws.Range("E" + ($startRow + $i)).Value = $periods[$i]
}
